$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 427, shifting existing rows 427:484 down to 428:485.
$ws.Rows.Item(427).Insert()

# Populate the newly inserted row 427 with the new weekly record.
$ws.Cells.Item(427, 1).Value = 5
$ws.Cells.Item(427, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(427, 3).Value = "Maule"
$ws.Cells.Item(427, 4).Value = Get-Date -Year 2023 -Month 2 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(427, 5).Value = 7
$ws.Cells.Item(427, 6).Value = 100114014
$ws.Cells.Item(427, 7).Value = "Betarraga"
$ws.Cells.Item(427, 8).Value = "Sin especificar"
$ws.Cells.Item(427, 9).Value = "Primera"
$ws.Cells.Item(427, 10).Value = 5000
$ws.Cells.Item(427, 11).Value = 600
$ws.Cells.Item(427, 12).Value = 600
$ws.Cells.Item(427, 13).Value = 600
$ws.Cells.Item(427, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(427, 15).Value = "Región del Maule"
$ws.Cells.Item(427, 16).Value = 120
$ws.Cells.Item(427, 17).Value = 5
$ws.Cells.Item(427, 18).Value = "Hortaliza"

# Match the date cell formatting used by the rest of column D.
$ws.Cells.Item(427, 4).NumberFormat = $ws.Cells.Item(428, 4).NumberFormat
